$d = $word.ActiveDocument

$d.Content.Find.Execute("54×44=2376", $true, $false, $false, $false, $false, $true, 1, $false, "26×38=988", 2) | Out-Null
$d.Content.Find.Execute("33×65=2145", $true, $false, $false, $false, $false, $true, 1, $false, "48×99=4752", 2) | Out-Null
$d.Content.Find.Execute("46×18=828", $true, $false, $false, $false, $false, $true, 1, $false, "32×98=3136", 2) | Out-Null
$d.Content.Find.Execute("70×55=3850", $true, $false, $false, $false, $false, $true, 1, $false, "54×97=5238", 2) | Out-Null
$d.Content.Find.Execute("92×13=1196", $true, $false, $false, $false, $false, $true, 1, $false, "90×89=8010", 2) | Out-Null
$d.Content.Find.Execute("76×18=1368", $true, $false, $false, $false, $false, $true, 1, $false, "49×39=1911", 2) | Out-Null
$d.Content.Find.Execute("83×23=1909", $true, $false, $false, $false, $false, $true, 1, $false, "64×71=4544", 2) | Out-Null
$d.Content.Find.Execute("82×49=4018", $true, $false, $false, $false, $false, $true, 1, $false, "56×71=3976", 2) | Out-Null
$d.Content.Find.Execute("71×58=4118", $true, $false, $false, $false, $false, $true, 1, $false, "48×27=1296", 2) | Out-Null
$d.Content.Find.Execute("89×14=1246", $true, $false, $false, $false, $false, $true, 1, $false, "17×78=1326", 2) | Out-Null
$d.Content.Find.Execute("89×35=3115", $true, $false, $false, $false, $false, $true, 1, $false, "64×95=6080", 2) | Out-Null
$d.Content.Find.Execute("73×77=5621", $true, $false, $false, $false, $false, $true, 1, $false, "80×14=1120", 2) | Out-Null
$d.Content.Find.Execute("15×50=750", $true, $false, $false, $false, $false, $true, 1, $false, "90×84=7560", 2) | Out-Null
$d.Content.Find.Execute("48×94=4512", $true, $false, $false, $false, $false, $true, 1, $false, "38×31=1178", 2) | Out-Null
$d.Content.Find.Execute("58×50=2900", $true, $false, $false, $false, $false, $true, 1, $false, "50×61=3050", 2) | Out-Null
$d.Content.Find.Execute("29×23=667", $true, $false, $false, $false, $false, $true, 1, $false, "75×90=6750", 2) | Out-Null
$d.Content.Find.Execute("64×80=5120", $true, $false, $false, $false, $false, $true, 1, $false, "67×69=4623", 2) | Out-Null
$d.Content.Find.Execute("23×64=1472", $true, $false, $false, $false, $false, $true, 1, $false, "12×27=324", 2) | Out-Null
$d.Content.Find.Execute("92×38=3496", $true, $false, $false, $false, $false, $true, 1, $false, "99×50=4950", 2) | Out-Null
$d.Content.Find.Execute("86×65=5590", $true, $false, $false, $false, $false, $true, 1, $false, "86×71=6106", 2) | Out-Null
$d.Content.Find.Execute("18×79=1422", $true, $false, $false, $false, $false, $true, 1, $false, "70×16=1120", 2) | Out-Null
$d.Content.Find.Execute("99×87=8613", $true, $false, $false, $false, $false, $true, 1, $false, "26×47=1222", 2) | Out-Null
$d.Content.Find.Execute("43×75=3225", $true, $false, $false, $false, $false, $true, 1, $false, "96×84=8064", 2) | Out-Null
$d.Content.Find.Execute("74×40=2960", $true, $false, $false, $false, $false, $true, 1, $false, "22×21=462", 2) | Out-Null
$d.Content.Find.Execute("38×74=2812", $true, $false, $false, $false, $false, $true, 1, $false, "61×63=3843", 2) | Out-Null
